# Add new "Transformation Energy" / "Final Energy" rows to the
# variable_definitions sheet (industry energy plus markdown annotations).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variable_definitions")

$newVars = @(
    "Transformation Energy|Coal",
    "Transformation Energy|Oil",
    "Transformation Energy|Gas",
    "Transformation Energy|Biomass",
    "Final Energy|Industry|Gas Feedstocks",
    "Final Energy|Industry|Oil Feedstocks"
)

$startRow = 143
for ($i = 0; $i -lt $newVars.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "NA"
    $ws.Cells.Item($r, 2).Value = "Energy"
    $ws.Cells.Item($r, 3).Value = "NA"
    $ws.Cells.Item($r, 4).Value = $newVars[$i]
    $ws.Cells.Item($r, 5).Value = "EJ/yr"
    $ws.Cells.Item($r, 6).Value = "created for LEEP report data"
}

# Mirror the author's final on-screen selection/scroll position.
$ws.Range("F137").Select()
